# Fix spelling/capitalization mistakes in the "Matriz de seguimiento" sheet
# and tidy up a couple of presentation details (selection, page orientation).
#
# NOTE: the order in which new cell text values are assigned matters because
# each newly-introduced distinct string is appended to the shared string
# table the first time it is written; the order below reproduces the shared
# string table layout of the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value  = "Inicio de sesión"
$ws.Range("A4").Value  = "Página inicial usuario"
$ws.Range("A6").Value  = "Registro de usuario"
$ws.Range("A8").Value  = "Forma de pago"
$ws.Range("A9").Value  = "Vender"
$ws.Range("A10").Value = "Nosotros"
$ws.Range("A12").Value = "Perfil de usuario"
$ws.Range("A13").Value = "Seccion de anime"
$ws.Range("A7").Value  = "Carro de compras"
$ws.Range("A11").Value = "Página de los productos"
$ws.Range("A15").Value = "Seccion de electrónica"
$ws.Range("F2").Value  = "Sin empezar"
$ws.Range("A3").Value  = "Página Inicial (Index)"

# --- Update the active selection --------------------------------------------

$ws.Range("A3:B3").Select()

# --- Set explicit page orientation (portrait) -------------------------------

$ws.PageSetup.Orientation = 1
